# The "missing_data.xlsx" table re-runs its error calculation / plotting
# pass: two rows of raw measurements (the "RM 232" row and the "SC 92" row)
# are dropped from the bottom block, the remaining "SC *" rows shift up by
# one, and the per-cell "is this value missing" mask is redrawn, filling
# some previously-blank cells and blanking a few previously-filled ones.
# Net effect: the used range shrinks from A1:F35 to A1:F33.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the last two rows of the original table so the sheet ends at F33.
$ws.Range("A34:F35").EntireRow.Delete()

# --- cells whose "missing" mask flips in the RM block (rows 2-23) -------
$ws.Range("C2").Value = 14.9
$ws.Range("D3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = -15.4
$ws.Range("D5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("C12").Value = 12.5
$ws.Range("C14").ClearContents()
$ws.Range("C20").Value = 12.5
$ws.Range("C21").Value = 12.7
$ws.Range("C22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("D23").Value = -13.9

# --- bottom block (rows 26-33): "RM 232" and "SC 92" rows removed, -------
# --- remaining SC rows shift up one position with a refreshed mask -------
$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = 10
$ws.Range("D27").ClearContents()
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = 17

$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").Value = -19.6
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44

$ws.Range("A29").Value = "SC 119"
$ws.Range("B29").Value = -19.5
$ws.Range("C29").Value = 11.2
$ws.Range("D29").Value = -13
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = 18.06

$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18

$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").ClearContents()
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53
